# Update the generation Date shown on the Metadata sheet.
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-06-04T08:55:54+00:00"

# Elements sheet: add a binding (Binding Strength / Binding Value Set) to the
# Extension.value[x] row (row 6), and widen the "Binding Value Set" column to
# fit the new long value-set URL.
$ws = $wb.Worksheets.Item("Elements")
$ws.Range("X6").Value = "required"
$ws.Range("Y6").Value = ""
$ws.Range("Z6").Value = "http://ltsi.univ-rennes.fr/ValueSet/siph-listeald-oncofair-valueset"
$ws.Columns.Item(26).ColumnWidth = 59
